$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Home win": append a new last row (row 3) with the new match
# ---------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")
$wsHome.Cells.Item(3,1).Value = "10-05-2025 12:00"
$wsHome.Cells.Item(3,2).Value = "CZECH-REPUBLIC"
$wsHome.Cells.Item(3,3).Value = "1. LIGA U19"
$wsHome.Cells.Item(3,4).Value = "Zbrojovka Brno U19 - Slovácko U19"
$wsHome.Cells.Item(3,5).Value = 73.3
$wsHome.Cells.Item(3,6).Value = 2.75

# ---------------------------------------------------------------
# Sheet "EV Home win": insert 4 new rows at their chronological spots
# ---------------------------------------------------------------
$wsEvHome = $wb.Worksheets.Item("EV Home win")

# New row 4: AUSTRALIA (pushes old rows 4-10 down to 5-11)
$wsEvHome.Rows.Item(4).Insert()
$wsEvHome.Rows.Item(4).ClearFormats()
$wsEvHome.Cells.Item(4,1).Value = "10-05-2025 12:00"
$wsEvHome.Cells.Item(4,2).Value = "AUSTRALIA"
$wsEvHome.Cells.Item(4,3).Value = "NEW SOUTH WALES NPL 2"
$wsEvHome.Cells.Item(4,4).Value = "Mounties Wanderers - Hills Brumbies"
$wsEvHome.Cells.Item(4,5).Value = 50
$wsEvHome.Cells.Item(4,6).Value = 2.9
$wsEvHome.Cells.Item(4,7).Value = 0.45

# New rows 7 & 8: CZECH-REPUBLIC matches (pushes old rows down by 2 more)
$wsEvHome.Rows.Item(7).Insert()
$wsEvHome.Rows.Item(7).ClearFormats()
$wsEvHome.Cells.Item(7,1).Value = "10-05-2025 12:00"
$wsEvHome.Cells.Item(7,2).Value = "CZECH-REPUBLIC"
$wsEvHome.Cells.Item(7,3).Value = "1. LIGA U19"
$wsEvHome.Cells.Item(7,4).Value = "Zbrojovka Brno U19 - Slovácko U19"
$wsEvHome.Cells.Item(7,5).Value = 73.3
$wsEvHome.Cells.Item(7,6).Value = 2.75
$wsEvHome.Cells.Item(7,7).Value = 1.02

$wsEvHome.Rows.Item(8).Insert()
$wsEvHome.Rows.Item(8).ClearFormats()
$wsEvHome.Cells.Item(8,1).Value = "10-05-2025 12:00"
$wsEvHome.Cells.Item(8,2).Value = "CZECH-REPUBLIC"
$wsEvHome.Cells.Item(8,3).Value = "1. LIGA U19"
$wsEvHome.Cells.Item(8,4).Value = "Sigma Olomouc U19 - Zlín U19"
$wsEvHome.Cells.Item(8,5).Value = 55.7
$wsEvHome.Cells.Item(8,6).Value = 2.1
$wsEvHome.Cells.Item(8,7).Value = 0.17

# New row 13: VIETNAM (pushes old last row down to 14)
$wsEvHome.Rows.Item(13).Insert()
$wsEvHome.Rows.Item(13).ClearFormats()
$wsEvHome.Cells.Item(13,1).Value = "10-05-2025 12:00"
$wsEvHome.Cells.Item(13,2).Value = "VIETNAM"
$wsEvHome.Cells.Item(13,3).Value = "V.LEAGUE 2"
$wsEvHome.Cells.Item(13,4).Value = "Dong Nai - Dong Thap"
$wsEvHome.Cells.Item(13,5).Value = 60
$wsEvHome.Cells.Item(13,6).Value = 2.2
$wsEvHome.Cells.Item(13,7).Value = 0.32

# ---------------------------------------------------------------
# Sheet "EV Away win": insert a new row 2 (AUSTRALIA) at the top
# ---------------------------------------------------------------
$wsEvAway = $wb.Worksheets.Item("EV Away win")
$wsEvAway.Rows.Item(2).Insert()
$wsEvAway.Rows.Item(2).ClearFormats()
$wsEvAway.Cells.Item(2,1).Value = "10-05-2025 12:00"
$wsEvAway.Cells.Item(2,2).Value = "AUSTRALIA"
$wsEvAway.Cells.Item(2,3).Value = "NEW SOUTH WALES NPL 2"
$wsEvAway.Cells.Item(2,4).Value = "SD Raiders - Blacktown Spartans"
$wsEvAway.Cells.Item(2,5).Value = 51
$wsEvAway.Cells.Item(2,6).Value = 3.85
$wsEvAway.Cells.Item(2,7).Value = 0.96

# ---------------------------------------------------------------
# Sheet "EV Btts": insert a new row 2 (ITALY) at the top
# ---------------------------------------------------------------
$wsEvBtts = $wb.Worksheets.Item("EV Btts")
$wsEvBtts.Rows.Item(2).Insert()
$wsEvBtts.Rows.Item(2).ClearFormats()
$wsEvBtts.Cells.Item(2,1).Value = "10-05-2025 12:00"
$wsEvBtts.Cells.Item(2,2).Value = "ITALY"
$wsEvBtts.Cells.Item(2,3).Value = "CAMPIONATO PRIMAVERA - 1"
$wsEvBtts.Cells.Item(2,4).Value = "Internazionale U20 - Sampdoria U20"
$wsEvBtts.Cells.Item(2,5).Value = 68
$wsEvBtts.Cells.Item(2,6).Value = 1.95
$wsEvBtts.Cells.Item(2,7).Value = 0.33
